# Update the forecast dates (column A) and scores (column C) for rows 2-31
# to reflect the refreshed scrapper output, per the commit "Scrapper ready
# for production". Signal text in column B is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    45940, 45943, 45944, 45945, 45946, 45947, 45950, 45951, 45952, 45953,
    45954, 45957, 45958, 45959, 45960, 45961, 45964, 45965, 45966, 45967,
    45968, 45971, 45972, 45973, 45974, 45975, 45978, 45979, 45980, 45981
)

$scores = @(
    0.2660646151183387, 0.2574717995054228, 0.2493119621046376, 0.2336090690140274,
    0.2349794051493647, 0.2330001286114695, 0.2344442068268122, 0.2359496702182244,
    0.2283995496571739, 0.2251089505700687, 0.2234114996938327, 0.2223929684345076,
    0.2291202871244187, 0.2319713214561409, 0.2338585535213417, 0.2315027899831548,
    0.2335542066693475, 0.235117210053213,  0.2351532779527054, 0.2429966940770025,
    0.2423682593055715, 0.2436330769571086, 0.2633896615943274, 0.2655242454178919,
    0.2664659839041819, 0.2609385276189906, 0.2586563026622665, 0.2612144140843821,
    0.2670336621723229, 0.2622518556396169
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 3).Value = $scores[$i]
}
